$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string value used by rows 20 and 21
$label = "Elkan Par Dyn"

# Row 20
$ws.Range("A20").Value = $label
$ws.Range("B20").Value = 2.163586
$ws.Range("C20").Value = 10.456343
$ws.Range("D20").Value = 30.44319
$ws.Range("E20").Value = 97.062709

# Row 21
$ws.Range("A21").Value = $label
$ws.Range("B21").Value = 2.432288
$ws.Range("C21").Value = 6.779244
$ws.Range("D21").Value = 30.706803
$ws.Range("E21").Value = 106.46536

# Row 22 - sums
$ws.Range("B22").Formula = "=SUM(B20:B21)"
$ws.Range("C22").Formula = "=SUM(C20:C21)"
$ws.Range("D22").Formula = "=SUM(D20:D21)"
$ws.Range("E22").Formula = "=SUM(E20:E21)"

# Row 23 - average (divide by 2)
$ws.Range("B23").Formula = "=B22/2"
$ws.Range("C23").Formula = "=C22/2"
$ws.Range("D23").Formula = "=D22/2"
$ws.Range("E23").Formula = "=E22/2"

# Update selection to I31
$ws.Range("I31").Select()
